# Update the workbook "Översikt HÖÖR" to reflect a new export run:
#  - Column C (Förändrad) moves from 2023-09-21 (45190) to 2023-09-23 (45192)
#    for every existing data row (2-172).
#  - A new case "A 44960-2023" is appended as row 173.
#  - The dimension / used range grows from A1:Y172 to A1:Y173.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" date for every existing data row (2..172) from
#    2023-09-21 to 2023-09-23 (serial 45190 -> 45192).
$ws.Range("C2:C172").Value = 45192

# Make sure row 172 carries an explicit row-height record, matching the
# rest of the sheet (rows 1-171 already have ht="15" customHeight="1").
$ws.Range("A172:Y172").RowHeight = 15

# 2) Append the new case as row 173.
$ws.Range("A173").Value = "A 44960-2023"

$ws.Range("B173").NumberFormat = "YYYY-MM-DD"
$ws.Range("B173").Value = 45190

$ws.Range("C173").NumberFormat = "YYYY-MM-DD"
$ws.Range("C173").Value = 45192

$ws.Range("D173").Value = "SKÅNE LÄN"
$ws.Range("E173").Value = "HÖÖR"

$ws.Range("G173").Value = 1.4
$ws.Range("H173").Value = 0
$ws.Range("I173").Value = 0
$ws.Range("J173").Value = 0
$ws.Range("K173").Value = 0
$ws.Range("L173").Value = 0
$ws.Range("M173").Value = 0
$ws.Range("N173").Value = 0
$ws.Range("O173").Value = 0
$ws.Range("P173").Value = 0
$ws.Range("Q173").Value = 0

# R column uses a wrap-text style throughout the sheet, even when blank.
$ws.Range("R173").WrapText = $true
